# Updates crypto price/volume figures (and re-syncs a handful of coin
# rows whose rank shuffled) to match the refreshed coinranking.com scrape,
# per the GitHub Actions commit "Updated symbol list on Mon Feb  6 19:53:08
# UTC 2023 with GitHub Actions".
#
# Price (D) and Volume(1h) (E) are plain-text columns in the source sheet
# (t="inlineStr"), e.g. "329.31" or "0.42%" - not real numbers/percentages.
# A bare assignment like $ws.Range("D2").Value = '329.22' would let Excel
# auto-detect the numeric/percent shape and silently convert the cell to a
# Number/Percentage, which changes both the stored value and its type. To
# keep these as text (matching the source workbook), values in columns D/E
# are written with a leading apostrophe, Excel's standard "force text"
# quote-prefix, exactly as if a user typed them in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Price/Volume(1h))
$ws.Range("D2").Value = '''329.22'
$ws.Range("E2").Value = '''0.28%'

# Row 3 (Price/Volume(1h))
$ws.Range("D3").Value = '''44.14'
$ws.Range("E3").Value = '''0.13%'

# Row 4 (Price/Volume(1h))
$ws.Range("D4").Value = '''5.571'
$ws.Range("E4").Value = '''2.21%'

# Row 5 (Price/Volume(1h))
$ws.Range("D5").Value = '''0.08097'
$ws.Range("E5").Value = '''0.24%'

# Row 6 (Price/Volume(1h))
$ws.Range("D6").Value = '''1.983'
$ws.Range("E6").Value = '''4.55%'

# Row 7 (Coin/Link/Price/Volume(1h))
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''4.319'
$ws.Range("E7").Value = '''0.46%'

# Row 8 (Coin/Link/Price/Volume(1h))
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.573'
$ws.Range("E8").Value = '''-6.25%'

# Row 9 (Coin/Link/Price/Volume(1h))
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9515'
$ws.Range("E9").Value = '''1.27%'

# Row 10 (Coin/Link/Price/Volume(1h))
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1170'
$ws.Range("E10").Value = '''-3.86%'

# Row 11 (Coin/Link/Price/Volume(1h))
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1848'
$ws.Range("E11").Value = '''-2.33%'

# Row 12 (Price/Volume(1h))
$ws.Range("D12").Value = '''0.09818'
$ws.Range("E12").Value = '''1.23%'

# Row 13 (Price/Volume(1h))
$ws.Range("D13").Value = '''0.04677'
$ws.Range("E13").Value = '''12.67%'

# Row 14 (Price/Volume(1h))
$ws.Range("D14").Value = '''0.1067'
$ws.Range("E14").Value = '''-0.27%'

# Row 15 (Price/Volume(1h))
$ws.Range("D15").Value = '''0.001284'
$ws.Range("E15").Value = '''0.83%'

# Row 16 (Price/Volume(1h))
$ws.Range("D16").Value = '''0.04213'
$ws.Range("E16").Value = '''-3.63%'

# Row 17 (Price/Volume(1h))
$ws.Range("D17").Value = '''0.005979'
$ws.Range("E17").Value = '''-2.00%'

# Row 18 (Price/Volume(1h))
$ws.Range("D18").Value = '''3.371'
$ws.Range("E18").Value = '''-5.62%'

# Row 19 (Coin/Link/Price/Volume(1h))
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3473'
$ws.Range("E19").Value = '''-0.72%'

# Row 20 (Coin/Link/Price/Volume(1h))
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '''10.21'
$ws.Range("E20").Value = '''19.94%'

# Row 21 (Volume(1h))
$ws.Range("E21").Value = '''4.41%'

# Row 22 (Price/Volume(1h))
$ws.Range("D22").Value = '''0.2506'
$ws.Range("E22").Value = '''0.48%'

# Row 23 (Price/Volume(1h))
$ws.Range("D23").Value = '''0.001249'
$ws.Range("E23").Value = '''1.08%'

# Row 24 (Price/Volume(1h))
$ws.Range("D24").Value = '''0.004323'
$ws.Range("E24").Value = '''0.73%'

# Row 25 (Price/Volume(1h))
$ws.Range("D25").Value = '''0.0001191'
$ws.Range("E25").Value = '''-3.51%'

# Row 26 (Volume(1h))
$ws.Range("E26").Value = '''-0.72%'

# Row 38 (Price/Volume(1h))
$ws.Range("D38").Value = '''0.02664'
$ws.Range("E38").Value = '''0.36%'

# Row 39 (Price/Volume(1h))
$ws.Range("D39").Value = '''0.05555'
$ws.Range("E39").Value = '''1.73%'

# Row 40 (Price/Volume(1h))
$ws.Range("D40").Value = '''0.007569'
$ws.Range("E40").Value = '''-1.51%'

# Row 41 (Volume(1h))
$ws.Range("E41").Value = '''1.32%'

# Row 42 (Volume(1h))
$ws.Range("E42").Value = '''-16.93%'

# Row 43 (Price/Volume(1h))
$ws.Range("D43").Value = '''0.002017'
$ws.Range("E43").Value = '''-5.14%'

# Row 44 (Price/Volume(1h))
$ws.Range("D44").Value = '''0.008873'
$ws.Range("E44").Value = '''-10.34%'

# Row 45 (Price/Volume(1h))
$ws.Range("D45").Value = '''0.00007219'
$ws.Range("E45").Value = '''1.32%'

# Row 46 (Price/Volume(1h))
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("E46").Value = '''-0.29%'

# Row 47 (Coin/Link/Price/Volume(1h))
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").Value = '''0.002271'
$ws.Range("E47").Value = '''-0.45%'

# Row 48 (Coin/Link/Price/Volume(1h))
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = '''0.004358'
$ws.Range("E48").Value = '''22.67%'

# Row 49 (Price/Volume(1h))
$ws.Range("D49").Value = '''0.00002101'
$ws.Range("E49").Value = '''-0.29%'

# Row 50 (Price/Volume(1h))
$ws.Range("D50").Value = '''0.0002001'
$ws.Range("E50").Value = '''-0.29%'
